$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - subject id headers
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - CON values
$ws.Range("B2").Value = 61.922882011607641
$ws.Range("C2").Value = 49.062128841075321
$ws.Range("D2").Value = 64.529399275861493
$ws.Range("E2").Value = 52.205685730939273

# Row 3 - STR values (C3 removed, D3 newly added)
$ws.Range("B3").Value = 64.28974475450562
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = 70.480641863599658
$ws.Range("E3").Value = 46.926664468428058

# Selection/view update reflecting the updated active range
$ws.Range("B1:E3").Select()
